# Refresh the cryptocurrency table (coin name, link, price, 1h volume %)
# with the latest scrape. Prices/links shift down a row where a new coin
# (OKB) was inserted into the ranking, and every price/percent column is
# re-synced to the newest values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; Coin = $null; Link = $null; Price = '28.116.98'; Volume = '  -4.51%  ' }
    @{ Row = 3; Coin = $null; Link = $null; Price = '1.830.49'; Volume = '  -3.31%  ' }
    @{ Row = 4; Coin = $null; Link = $null; Price = '0.9998'; Volume = '  -0.37%  ' }
    @{ Row = 5; Coin = $null; Link = $null; Price = '329.10'; Volume = '  -2.74%  ' }
    @{ Row = 6; Coin = $null; Link = $null; Price = $null; Volume = '  -0.36%  ' }
    @{ Row = 7; Coin = $null; Link = $null; Price = '0.4647'; Volume = '  -2.21%  ' }
    @{ Row = 8; Coin = $null; Link = $null; Price = $null; Volume = '  -3.25%  ' }
    @{ Row = 9; Coin = 'OKB'; Link = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; Price = '46.23'; Volume = $null }
    @{ Row = 10; Coin = 'Dogecoin'; Link = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'; Price = '0.07875'; Volume = '  -1.75%  ' }
    @{ Row = 11; Coin = 'Polygon'; Link = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'; Price = '0.9599'; Volume = '  -2.88%  ' }
    @{ Row = 12; Coin = 'Solana'; Link = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; Price = '22.00'; Volume = '  -4.80%  ' }
    @{ Row = 13; Coin = 'WrappedEther'; Link = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; Price = '1.880.81'; Volume = '  -2.19%  ' }
    @{ Row = 14; Coin = 'Polkadot'; Link = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'; Price = '5.664'; Volume = '  -4.35%  ' }
    @{ Row = 15; Coin = 'Chainlink'; Link = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'; Price = '6.895'; Volume = '  -2.50%  ' }
    @{ Row = 16; Coin = 'TRON'; Link = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; Price = '0.06858'; Volume = '  +0.79%  ' }
    @{ Row = 17; Coin = 'BinanceUSD'; Link = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'; Price = '0.9994'; Volume = '  -0.56%  ' }
    @{ Row = 18; Coin = 'Litecoin'; Link = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; Price = '86.73'; Volume = '  -2.52%  ' }
    @{ Row = 19; Coin = 'ShibaInu'; Link = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; Price = '0.000009986'; Volume = '  -1.91%  ' }
    @{ Row = 20; Coin = 'Avalanche'; Link = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'; Price = '16.65'; Volume = '  -3.69%  ' }
    @{ Row = 21; Coin = 'Dai'; Link = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; Price = '1.001'; Volume = '  -0.18%  ' }
    @{ Row = 22; Coin = 'WrappedBTC'; Link = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'; Price = '28.145.85'; Volume = '  -4.47%  ' }
    @{ Row = 23; Coin = 'Uniswap'; Link = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'; Price = '5.331'; Volume = '  -3.20%  ' }
    @{ Row = 24; Coin = 'Cosmos'; Link = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; Price = '11.01'; Volume = '  -5.38%  ' }
    @{ Row = 25; Coin = 'Toncoin'; Link = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; Price = '2.094'; Volume = '  -2.59%  ' }
    @{ Row = 26; Coin = 'WrappedliquidstakedEther2.0'; Link = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; Price = '2.075.53'; Volume = '  -3.48%  ' }
    @{ Row = 27; Coin = 'Monero'; Link = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; Price = '152.42'; Volume = '  -3.07%  ' }
    @{ Row = 28; Coin = 'EthereumClassic'; Link = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; Price = '19.23'; Volume = '  -2.00%  ' }
    @{ Row = 29; Coin = 'InternetComputer(DFINITY)'; Link = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; Price = '5.788'; Volume = '  -10.52%  ' }
    @{ Row = 30; Coin = 'LidoDAOToken'; Link = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'; Price = '1.975'; Volume = '  -3.52%  ' }
    @{ Row = 31; Coin = 'BitcoinCash'; Link = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; Price = '117.07'; Volume = '  -1.54%  ' }
    @{ Row = 32; Coin = 'ImmutableX'; Link = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; Price = '0.9386'; Volume = '  -5.61%  ' }
    @{ Row = 33; Coin = 'Stellar'; Link = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; Price = '0.09266'; Volume = '  -2.89%  ' }
    @{ Row = 34; Coin = 'Filecoin'; Link = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; Price = '5.305'; Volume = '  -2.96%  ' }
    @{ Row = 35; Coin = 'ARBITRUM'; Link = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; Price = '1.320'; Volume = '  -4.57%  ' }
    @{ Row = 36; Coin = 'HuobiToken'; Link = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'; Price = '3.348'; Volume = '  -5.08%  ' }
    @{ Row = 37; Coin = 'Hedera'; Link = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; Price = '0.05953'; Volume = '  -7.18%  ' }
    @{ Row = 38; Coin = 'VeChain'; Link = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; Price = '0.02150'; Volume = '  -3.96%  ' }
    @{ Row = 39; Coin = 'TrustWalletToken'; Link = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; Price = '1.147'; Volume = '  -4.38%  ' }
    @{ Row = 40; Coin = 'Frax'; Link = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'; Price = '0.9994'; Volume = '  -0.34%  ' }
    @{ Row = 41; Coin = 'FraxShare'; Link = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'; Price = '7.669'; Volume = '  -0.98%  ' }
    @{ Row = 42; Coin = 'TheSandbox'; Link = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; Price = '0.5598'; Volume = '  -3.86%  ' }
    @{ Row = 43; Coin = 'Aptos'; Link = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Price = '9.936'; Volume = '  -5.58%  ' }
    @{ Row = 44; Coin = 'Algorand'; Link = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; Price = '0.1773'; Volume = '  -2.50%  ' }
    @{ Row = 45; Coin = 'WEMIXToken'; Link = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'; Price = '1.226'; Volume = '  -2.99%  ' }
    @{ Row = 46; Coin = 'RenderToken'; Link = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'; Price = '2.219'; Volume = '  -8.56%  ' }
    @{ Row = 47; Coin = 'EnergySwap'; Link = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; Price = '11.65'; Volume = '  -4.25%  ' }
    @{ Row = 48; Coin = 'Decentraland'; Link = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'; Price = '0.5271'; Volume = '  -4.01%  ' }
    @{ Row = 49; Coin = 'Cronos'; Link = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; Price = '0.07035'; Volume = '  -4.09%  ' }
    @{ Row = 50; Coin = 'NEARProtocol'; Link = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; Price = '1.835'; Volume = '  -5.88%  ' }
    @{ Row = 51; Coin = 'Quant'; Link = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'; Price = '111.68'; Volume = '  -3.99%  ' }
)

foreach ($r in $rows) {
    if ($r.Coin)   { $ws.Cells.Item($r.Row, 2).Value = $r.Coin }
    if ($r.Link)   { $ws.Cells.Item($r.Row, 3).Value = $r.Link }
    if ($r.Price)  {
        $priceCell = $ws.Cells.Item($r.Row, 4)
        $priceCell.NumberFormat = "@"   # keep dotted/zero-padded price strings literal
        $priceCell.Value = $r.Price
    }
    if ($r.Volume) { $ws.Cells.Item($r.Row, 5).Value = $r.Volume }
}
